$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 360; this shifts the existing rows
# 360-377 down to 362-379 (same as rows shifting down in the diff).
$ws.Rows("360:361").Insert()

# New row 360: weekly price entry (Primera)
$ws.Range("A360").Value = 11
$ws.Range("B360").Value = "Vega Monumental Concepción"
$ws.Range("C360").Value = "Bíobío"
$ws.Range("D360").Value = 44826
$ws.Range("E360").Value = 8
$ws.Range("F360").Value = 100114014
$ws.Range("G360").Value = "Betarraga"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 600
$ws.Range("K360").Value = 700
$ws.Range("L360").Value = 800
$ws.Range("M360").Value = 750
$ws.Range("N360").Value = "$/paquete 5 unidades"
$ws.Range("O360").Value = "Región Metropolitana"
$ws.Range("P360").Value = 150
$ws.Range("Q360").Value = 5
$ws.Range("R360").Value = "Hortaliza"

# New row 361: weekly price entry (Segunda)
$ws.Range("A361").Value = 11
$ws.Range("B361").Value = "Vega Monumental Concepción"
$ws.Range("C361").Value = "Bíobío"
$ws.Range("D361").Value = 44826
$ws.Range("E361").Value = 8
$ws.Range("F361").Value = 100114014
$ws.Range("G361").Value = "Betarraga"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Segunda"
$ws.Range("J361").Value = 300
$ws.Range("K361").Value = 600
$ws.Range("L361").Value = 600
$ws.Range("M361").Value = 600
$ws.Range("N361").Value = "$/paquete 5 unidades"
$ws.Range("O361").Value = "Región Metropolitana"
$ws.Range("P361").Value = 120
$ws.Range("Q361").Value = 5
$ws.Range("R361").Value = "Hortaliza"
